$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix typo in the "4. ..." row: "แยกภาพ" -> "ฝึกแยกภาพ", and give it the same
#    highlighted "headline" formatting as the other numbered steps (it was plain before).
$ws.Range("B17").Value = "4.ทดลองฝึกแยกภาพด้วย TensorFlow ใช้ dataSet : MNIST"
$ws.Range("B15").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2) Append the new task rows under row 18.
$ws.Range("B19").Value = "5.ทดลองฝึกแยกตัวอักษรด้วย TensorFlow ใช้ dataSet : IMDB"
$ws.Range("B20").Value = "Path : Research/lab/basic_classification_text.py"
$ws.Range("B21").Value = "6.ทดลอง Plot graph"
$ws.Range("B22").Value = "Path : Research/lab/tutorial_plotGraph.py"
$ws.Range("B23").Value = "7.จ่ายงาน"
$ws.Range("B24").Value = "7.1)Research"
$ws.Range("B25").Value = "7.1.1)tensorflow"
$ws.Range("B26").Value = "7.1.2)numpy"
$ws.Range("B27").Value = "7.1.3)keras"

# 3) The "headline" rows (3./4./5./6. ...) use the highlighted fill style that
#    already exists on B15/B17 - copy that formatting onto the new headline rows.
$ws.Range("B15").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B15").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4) Switch the workbook's base font from Calibri to Tahoma.
$wb.Styles.Item("Normal").Font.Name = "Tahoma"

# 5) Column widths were nudged slightly as part of this edit.
$ws.Columns.Item(1).ColumnWidth = 12.25
$ws.Columns.Item(2).ColumnWidth = 73.125
$ws.Columns.Item(3).ColumnWidth = 36.75

# 6) Update the view: scroll so row 11 is at the top and select the last entry.
$ws.Range("B11").Select()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B27").Select()
